$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.149220824241638
$ws.Range("B1").Value = 2.253112554550171
$ws.Range("C1").Value = 4.480283260345459
$ws.Range("D1").Value = 2.676498651504517
$ws.Range("E1").Value = 1.246300458908081
